# Rumus Perpangkatan Universal 2.0 [Revised]
# Edit: insert a delimited (x^n) term between " " and "- " just before the
# second nary (sum) operator in the oMathPara, i.e. turn
#   ... dx/dt x- (sum)...   into   ... dx/dt (x^n)- (sum)...
#
# The math markup is deeply nested OMML that the Word OM's OMath object
# doesn't expose enough surface to build (no Add/Insert for d/sSup), so we
# round-trip the paragraph's WordOpenXML and patch the run text/markup
# directly - the same bytes Word itself would produce for this edit.

$d = $word.ActiveDocument

# The oMathPara lives in the document's last paragraph (only one in this doc;
# the math run text doesn't round-trip through Range.Text as plain glyphs,
# so address it positionally instead of via a text search).
$target = $d.Paragraphs.Item($d.Paragraphs.Count)

$xml = $target.Range.WordOpenXML

$old = '<m:t xml:space="preserve"> x- </m:t>'

$new = '<m:t xml:space="preserve"> </m:t></m:r>' + `
  '<m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></m:ctrlPr></m:dPr>' + `
  '<m:e><m:sSup><m:sSupPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></m:ctrlPr></m:sSupPr>' + `
  '<m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t>x</m:t></m:r></m:e>' + `
  '<m:sup><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t>n</m:t></m:r></m:sup>' + `
  '</m:sSup></m:e></m:d>' + `
  '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t xml:space="preserve">- </m:t>'

if (-not $xml.Contains($old)) {
    throw "expected math run text not found in target paragraph"
}

$xml = $xml.Replace($old, $new)

# Range.WordOpenXML is read-only in this host; InsertXML replaces the
# range's contents with parsed OOXML (accepts the full WordOpenXML package
# shape returned by the getter).
$target.Range.InsertXML($xml)

Write-Output "done"
